$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 6 "comment" cells below each had a full stop (.) added to the end of
# their text (and, for one, a trailing space trimmed). Editing the cell text
# in Excel causes the shared-string table to be rewritten: the edited
# strings move to the end of the table while everything else shifts down to
# fill the gaps - this is exactly what the target diff shows.

$ws.Range("C3").Value = "Emissions decline by over 80% as mainly off road liquid fuels are displaced. Some diesel still remains from hard to decarbonise areas such as high utilisation tractors and commercial fishing vessels."
$ws.Range("C11").Value = "Emissions decline in both scenarios - but some natural gas remains in water and space heating, having a larger impact in Kea due to the higher GDP growth."
$ws.Range("C22").Value = "Coal and Natural Gas decrease significantly as new wind and hydro generation gets developed. Natural Gas is supplied during winter months due to the limited amount of other economically achievable options."
$ws.Range("C24").Value = "Large decarbonisation occurs in the industrial sector with most of the remaining fossil fuels in hard to abate sectors. Note - Only energy related emissions are included in the TIMES-NZ model. Emissions from Feedstock are not expressed."
$ws.Range("C26").Value = "Diesel prolongs in both scenarios as alternative technologies are not yet available. There is 4 PJ of biofuel in Kea further reducing diesel emissions."
$ws.Range("C32").Value = "Diesel prolongs in both scenarios as alternative technologies are not yet available. There is 2.3 PJ of biofuel in Kea further reducing diesel emissions."

# Update the active selection / scroll position to match the saved view:
# scrolled down so row 29 is the top-left visible row, with C48 selected.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 29
$aw.ScrollColumn = 1
$ws.Range("C48").Select()
